# Checkpoint before follow-up message
# Updates the "Active Signals", "Summary Dashboard" and "Signal History"
# sheets of genx_signals.xlsx with the latest signal snapshot.

$wb = $excel.ActiveWorkbook

$wsActive  = $wb.Worksheets.Item("Active Signals")
$wsSummary = $wb.Worksheets.Item("Summary Dashboard")
$wsHistory = $wb.Worksheets.Item("Signal History")

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Helper: write a value into a cell while forcing it to stay plain TEXT
# (Excel would otherwise "smart parse" numeric-looking / percent-looking
# strings into numbers). After writing, the cell's visual style is restored
# by pasting the formatting (only) from a nearby cell that already carries
# the desired style, so borders/fills/fonts are preserved exactly.
# ---------------------------------------------------------------------------
function Set-TextValue($ws, $cellRef, $value, $styleSourceRef) {
    $target = $ws.Range($cellRef)
    $target.NumberFormat = "@"
    $target.Value = $value
    $ws.Range($styleSourceRef).Copy()
    $target.PasteSpecial($xlPasteFormats)
}

# ===========================================================================
# Sheet: Active Signals
# ===========================================================================
# columns: A Timestamp, B Symbol, C Signal, D Entry, E SL, F TP, G Lots,
#          H Confidence (text, e.g. "77.0%"), I R:R, J Status

$activeRows = @(
    @{ Row=2; Timestamp="2025-07-28 19:44"; Symbol="AUDUSD"; Signal="BUY";  Entry=0.6569199999999999; SL=0.65464; TP=0.66323; Lots=0.05; Confidence="77.0%"; RR=2.77 },
    @{ Row=3; Timestamp="2025-07-28 19:29"; Symbol="USDCAD"; Signal="BUY";  Entry=1.36116;            SL=1.35649; TP=1.36818; Lots=0.09; Confidence="68.0%"; RR=1.51 },
    @{ Row=4; Timestamp="2025-07-28 19:24"; Symbol="GBPUSD"; Signal="SELL"; Entry=1.27154;            SL=1.27653; TP=1.26573; Lots=0.08; Confidence="91.0%"; RR=1.16 },
    @{ Row=5; Timestamp="2025-07-28 19:40"; Symbol="AUDUSD"; Signal="BUY";  Entry=0.65813;            SL=0.65364; TP=0.6622400000000001; Lots=0.07000000000000001; Confidence="65.0%"; RR=0.92 }
)

foreach ($r in $activeRows) {
    $row = $r.Row

    $wsActive.Cells.Item($row, 1).Value = $r.Timestamp
    $wsActive.Cells.Item($row, 2).Value = $r.Symbol

    # Column C ("Signal") carries a conditional fill (green for BUY, pink for
    # SELL); restyle it from a row that already has the matching look.
    $cCell = $wsActive.Cells.Item($row, 3)
    $cCell.Value = $r.Signal
    if ($r.Signal -eq "BUY") {
        $wsActive.Range("C2").Copy()
    } else {
        $wsActive.Range("C4").Copy()
    }
    $cCell.PasteSpecial($xlPasteFormats)

    $wsActive.Cells.Item($row, 4).Value = $r.Entry
    $wsActive.Cells.Item($row, 5).Value = $r.SL
    $wsActive.Cells.Item($row, 6).Value = $r.TP
    $wsActive.Cells.Item($row, 7).Value = $r.Lots

    # Column H ("Confidence") is plain text like "77.0%" - keep it text.
    Set-TextValue $wsActive ("H" + $row) $r.Confidence "G4"

    $wsActive.Cells.Item($row, 9).Value = $r.RR
    # Column J ("Status") is unchanged ("Active") for every row - leave as is.
}

# ===========================================================================
# Sheet: Summary Dashboard
# ===========================================================================

$wsSummary.Range("B5").Value = 10
$wsSummary.Range("B6").Value = 5

Set-TextValue $wsSummary "B7" "79.6%" "B3"
Set-TextValue $wsSummary "B8" "1.99" "B3"
Set-TextValue $wsSummary "B9" "2025-07-28 19:25:28" "B3"

# ===========================================================================
# Sheet: Signal History
# ===========================================================================
# columns: A Timestamp, B Symbol, C Signal, D Entry, E SL, F TP, G Lots,
#          H Confidence (plain number 0-1), I R:R, J Status

$historyRows = @(
    @{ Row=2;  Timestamp="2025-07-28 19:44"; Symbol="AUDUSD"; Signal="BUY";  Entry=0.6569199999999999; SL=0.65464; TP=0.66323; Lots=0.05; Confidence=0.77; RR=2.77; Status="Active" },
    @{ Row=3;  Timestamp="2025-07-28 19:29"; Symbol="USDCAD"; Signal="BUY";  Entry=1.36116;            SL=1.35649; TP=1.36818; Lots=0.09; Confidence=0.68; RR=1.51; Status="Active" },
    @{ Row=4;  Timestamp="2025-07-28 19:38"; Symbol="USDCAD"; Signal="BUY";  Entry=1.36194;            SL=1.35947; TP=1.3684;  Lots=0.1;  Confidence=0.77; RR=2.61; Status="Pending" },
    @{ Row=5;  Timestamp="2025-07-28 19:24"; Symbol="GBPUSD"; Signal="SELL"; Entry=1.27154;            SL=1.27653; TP=1.26573; Lots=0.08; Confidence=0.91; RR=1.16; Status="Active" },
    @{ Row=6;  Timestamp="2025-07-28 19:34"; Symbol="USDJPY"; Signal="BUY";  Entry=150.00672;          SL=149.69719; TP=150.70286; Lots=0.09; Confidence=0.79; RR=2.25; Status="Filled" },
    @{ Row=7;  Timestamp="2025-07-28 19:02"; Symbol="GBPUSD"; Signal="BUY";  Entry=1.27111;            SL=1.26667; TP=1.28048; Lots=0.08; Confidence=0.82; RR=2.11; Status="Pending" },
    @{ Row=8;  Timestamp="2025-07-28 19:11"; Symbol="NZDUSD"; Signal="BUY";  Entry=0.58896;            SL=0.5866;  TP=0.59733; Lots=0.04; Confidence=0.82; RR=3.56; Status="Filled" },
    @{ Row=9;  Timestamp="2025-07-28 19:03"; Symbol="EURUSD"; Signal="SELL"; Entry=1.10897;            SL=1.11237; TP=1.10414; Lots=0.09; Confidence=0.74; RR=1.42; Status="Filled" },
    @{ Row=10; Timestamp="2025-07-28 19:02"; Symbol="AUDUSD"; Signal="BUY";  Entry=0.65579;            SL=0.65228; TP=0.6625799999999999; Lots=0.03; Confidence=0.79; RR=1.94; Status="Filled" },
    @{ Row=11; Timestamp="2025-07-28 19:21"; Symbol="AUDUSD"; Signal="SELL"; Entry=0.65563;            SL=0.6601;  TP=0.65027; Lots=0.02; Confidence=0.89; RR=1.2;  Status="Filled" },
    @{ Row=12; Timestamp="2025-07-28 19:41"; Symbol="NZDUSD"; Signal="SELL"; Entry=0.59095;            SL=0.59576; TP=0.58178; Lots=0.02; Confidence=0.91; RR=1.91; Status="Pending" },
    @{ Row=13; Timestamp="2025-07-28 19:40"; Symbol="AUDUSD"; Signal="BUY";  Entry=0.65813;            SL=0.65364; TP=0.6622400000000001; Lots=0.07000000000000001; Confidence=0.65; RR=0.92; Status="Active" },
    @{ Row=14; Timestamp="2025-07-28 19:24"; Symbol="GBPUSD"; Signal="BUY";  Entry=1.27136;            SL=1.26819; TP=1.27863; Lots=0.06; Confidence=0.82; RR=2.3;  Status="Filled" },
    @{ Row=15; Timestamp="2025-07-28 19:24"; Symbol="USDCHF"; Signal="BUY";  Entry=0.87943;            SL=0.87619; TP=0.88415; Lots=0.07000000000000001; Confidence=0.88; RR=1.45; Status="Filled" },
    @{ Row=16; Timestamp="2025-07-28 19:21"; Symbol="GBPUSD"; Signal="SELL"; Entry=1.26355;            SL=1.26715; TP=1.25387; Lots=0.07000000000000001; Confidence=0.7;  RR=2.69; Status="Pending" }
)

foreach ($r in $historyRows) {
    $row = $r.Row
    $wsHistory.Cells.Item($row, 1).Value = $r.Timestamp
    $wsHistory.Cells.Item($row, 2).Value = $r.Symbol
    $wsHistory.Cells.Item($row, 3).Value = $r.Signal
    $wsHistory.Cells.Item($row, 4).Value = $r.Entry
    $wsHistory.Cells.Item($row, 5).Value = $r.SL
    $wsHistory.Cells.Item($row, 6).Value = $r.TP
    $wsHistory.Cells.Item($row, 7).Value = $r.Lots
    $wsHistory.Cells.Item($row, 8).Value = $r.Confidence
    $wsHistory.Cells.Item($row, 9).Value = $r.RR
    $wsHistory.Cells.Item($row, 10).Value = $r.Status
}
